$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the numeric-looking Price/Volume(1h) text columns from Excel automatic
# "looks like a number" conversion while we write the new text values.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '25.847.86'
$ws.Range('E2').Value = '  -3.55%  '
$ws.Range('D3').Value = '1.822.05'
$ws.Range('E3').Value = '  -2.58%  '
$ws.Range('D4').Value = '0.9956'
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').Value = '278.79'
$ws.Range('E5').Value = '  -7.20%  '
$ws.Range('D6').Value = '0.9965'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').Value = '0.5101'
$ws.Range('E7').Value = '  -4.27%  '
$ws.Range('D8').Value = '0.3469'
$ws.Range('E8').Value = '  -7.02%  '
$ws.Range('D9').Value = '44.45'
$ws.Range('E9').Value = '  -1.89%  '
$ws.Range('D10').Value = '0.06788'
$ws.Range('E10').Value = '  -5.01%  '
$ws.Range('D11').Value = '19.87'
$ws.Range('E11').Value = '  -7.74%  '
$ws.Range('D12').Value = '0.8082'
$ws.Range('E12').Value = '  -8.99%  '
$ws.Range('D13').Value = '0.07814'
$ws.Range('E13').Value = '  -3.99%  '
$ws.Range('D14').Value = '1.814.39'
$ws.Range('E14').Value = '  -3.16%  '
$ws.Range('D15').Value = '5.074'
$ws.Range('E15').Value = '  -4.13%  '
$ws.Range('D16').Value = '88.10'
$ws.Range('E16').Value = '  -4.85%  '
$ws.Range('D17').Value = '0.9941'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').Value = '14.13'
$ws.Range('E18').Value = '  -5.03%  '
$ws.Range('D19').Value = '0.000008054'
$ws.Range('E19').Value = '  -5.23%  '
$ws.Range('D20').Value = '0.9973'
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('D21').Value = '25.886.79'
$ws.Range('E21').Value = '  -3.48%  '
$ws.Range('D22').Value = '4.743'
$ws.Range('E22').Value = '  -4.63%  '
$ws.Range('D23').Value = '9.978'
$ws.Range('E23').Value = '  -6.22%  '
$ws.Range('D24').Value = '6.144'
$ws.Range('E24').Value = '  -3.74%  '
$ws.Range('D25').Value = '2.322'
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('D26').Value = '142.43'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('D27').Value = '1.659'
$ws.Range('E27').Value = '  -4.70%  '
$ws.Range('D28').Value = '17.16'
$ws.Range('E28').Value = '  -4.70%  '
$ws.Range('D29').Value = '109.29'
$ws.Range('E29').Value = '  -3.83%  '
$ws.Range('D30').Value = '4.332'
$ws.Range('E30').Value = '  -7.86%  '
$ws.Range('D31').Value = '4.286'
$ws.Range('E31').Value = '  -7.38%  '
$ws.Range('D32').Value = '0.08765'
$ws.Range('E32').Value = '  -3.73%  '
$ws.Range('D33').Value = '0.04855'
$ws.Range('E33').Value = '  -3.19%  '
$ws.Range('D34').Value = '1.164'
$ws.Range('E34').Value = '  -0.70%  '
$ws.Range('D35').Value = '0.7286'
$ws.Range('E35').Value = '  -9.95%  '
$ws.Range('D36').Value = '2.843'
$ws.Range('E36').Value = '  -3.52%  '
$ws.Range('D37').Value = '3.155'
$ws.Range('E37').Value = '  -0.91%  '
$ws.Range('D38').Value = '2.403'
$ws.Range('E38').Value = '  -9.36%  '
$ws.Range('D39').Value = '0.01845'
$ws.Range('E39').Value = '  -5.01%  '
$ws.Range('D40').Value = '0.5149'
$ws.Range('E40').Value = '  -15.82%  '
$ws.Range('D41').Value = '0.9484'
$ws.Range('E41').Value = '  -11.14%  '
$ws.Range('D42').Value = '116.72'
$ws.Range('E42').Value = '  +0.67%  '
$ws.Range('D43').Value = '6.202'
$ws.Range('E43').Value = '  -4.26%  '
$ws.Range('D44').Value = '7.969'
$ws.Range('E44').Value = '  -9.02%  '
$ws.Range('D45').Value = '0.9951'
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.4490'
$ws.Range('E46').Value = '  -15.48%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1358'
$ws.Range('E47').Value = '  -8.82%  '
$ws.Range('D48').Value = '9.252'
$ws.Range('E48').Value = '  -7.16%  '
$ws.Range('D49').Value = '36.20'
$ws.Range('E49').Value = '  -2.89%  '
$ws.Range('D50').Value = '0.05920'
$ws.Range('E50').Value = '  -2.29%  '
$ws.Range('D51').Value = '1.492'
$ws.Range('E51').Value = '  -9.86%  '

# Restore the plain Normal style so no extra number formatting is left behind.
$ws.Range("D2:E51").Style = "Normal"
